# Applies the football_weather data-refresh edit (FBS + Other sheets)
$wb = $excel.ActiveWorkbook

$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# --- FBS sheet: weather metric refresh ---
$wsFBS.Range("AB2").Value = 10
$wsFBS.Range("AF2").Value = 0.5
$wsFBS.Range("O2").Value = 36.85999999999999
$wsFBS.Range("P2").Value = 17.7
$wsFBS.Range("R2").Value = 0.4
$wsFBS.Range("U2").Value = 6
$wsFBS.Range("M3").Value = "SW"
$wsFBS.Range("N3").Value = "SW"
$wsFBS.Range("O3").Value = 36.02
$wsFBS.Range("Q3").Value = "SW"
$wsFBS.Range("M4").Value = "W"
$wsFBS.Range("N4").Value = "WNW"
$wsFBS.Range("O4").Value = 47.63
$wsFBS.Range("P4").Value = 11.5
$wsFBS.Range("Q4").Value = "WNW"
$wsFBS.Range("R4").Value = 0.1
$wsFBS.Range("S4").Value = 0
$wsFBS.Range("U4").Value = 4.6
$wsFBS.Range("M5").Value = "NW"
$wsFBS.Range("N5").Value = "WNW"
$wsFBS.Range("O5").Value = 43.52
$wsFBS.Range("P5").Value = 3.1
$wsFBS.Range("Q5").Value = "NW"
$wsFBS.Range("U5").Value = -2
$wsFBS.Range("M6").Value = "SE"
$wsFBS.Range("N6").Value = "SE"
$wsFBS.Range("Q6").Value = "SE"
$wsFBS.Range("N7").Value = "ESE"
$wsFBS.Range("Q7").Value = "ESE"
$wsFBS.Range("M8").Value = "WNW"
$wsFBS.Range("N8").Value = "W"
$wsFBS.Range("Q8").Value = "WNW"
$wsFBS.Range("N9").Value = "SE"
$wsFBS.Range("M17").Value = "W"
$wsFBS.Range("N17").Value = "WNW"
$wsFBS.Range("Q17").Value = "WNW"
$wsFBS.Range("M18").Value = "NNW"
$wsFBS.Range("N18").Value = "NW"
$wsFBS.Range("Q18").Value = "NNW"
$wsFBS.Range("N19").Value = "NNE"
$wsFBS.Range("M21").Value = "NW"
$wsFBS.Range("N21").Value = "NW"
$wsFBS.Range("Q21").Value = "NW"
$wsFBS.Range("M22").Value = "NW"
$wsFBS.Range("N22").Value = "NNW"
$wsFBS.Range("Q22").Value = "NNW"
$wsFBS.Range("M23").Value = "SE"
$wsFBS.Range("N23").Value = "SE"
$wsFBS.Range("Q23").Value = "SE"
$wsFBS.Range("M24").Value = "W"
$wsFBS.Range("N24").Value = "WNW"
$wsFBS.Range("Q24").Value = "WNW"
$wsFBS.Range("M25").Value = "W"
$wsFBS.Range("N25").Value = "W"
$wsFBS.Range("Q25").Value = "W"
$wsFBS.Range("M26").Value = "WSW"
$wsFBS.Range("N26").Value = "WSW"
$wsFBS.Range("Q26").Value = "WSW"
$wsFBS.Range("M28").Value = "NW"
$wsFBS.Range("N28").Value = "WNW"
$wsFBS.Range("Q28").Value = "WNW"
$wsFBS.Range("N29").Value = "NNE"
$wsFBS.Range("Q29").Value = "E"
$wsFBS.Range("M30").Value = "NW"
$wsFBS.Range("Q30").Value = "E"
$wsFBS.Range("N31").Value = "NNE"
$wsFBS.Range("Q31").Value = "E"
$wsFBS.Range("M33").Value = "NW"
$wsFBS.Range("N33").Value = "WNW"
$wsFBS.Range("Q33").Value = "WNW"
$wsFBS.Range("M34").Value = "W"
$wsFBS.Range("N34").Value = "NW"
$wsFBS.Range("Q34").Value = "SSW"

# --- Other sheet: weather metric refresh ---
$wsOther.Range("O2").Value = "NW"
$wsOther.Range("P2").Value = "NW"
$wsOther.Range("S2").Value = "NW"
$wsOther.Range("S4").Value = "NE"
$wsOther.Range("O5").Value = "SE"
$wsOther.Range("P5").Value = "WSW"
$wsOther.Range("S5").Value = "SE"

# --- Timestamp column refresh (shared across all FBS data rows) ---
$wsFBS.Range("AK2:AK34").Value = "2024-12-30T10:01:22.621044"
